# rerun module w updated plp
# Updates recalculated p-values in a handful of by-histology result sheets.

$wb = $excel.ActiveWorkbook

# --- Low-grade glioma ---
$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Range("C7").Value = 0.658631193852452
$ws.Range("C9").Value = 0.0809511288726465

# --- Medulloblastoma ---
$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Range("C7").Value = 0.0112805347433453

# --- Mixed neuronal-glial tumor ---
$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Range("C7").Value = 0.304722764052948
$ws.Range("C8").Value = 0.490780770910609

# --- ATRT ---
$ws = $wb.Worksheets.Item("ATRT")
$ws.Range("C2").Value = 0.696434129477608
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 0.705442383703247
$ws.Range("C5").Value = 0.490390189520623
$ws.Range("C6").Value = 0.447863247863246
$ws.Range("C7").Value = 0.260831159686995
$ws.Range("C8").Value = 0.285246835544318
